$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the existing 2026/02/21 block (row 854),
# shifting rows 855-896 down to 857-898.
$ws.Rows("855:856").Insert()

# New row 855: 2026/02/21 (Sat), 22:00, ranking 201
$ws.Range("A855").Value = "'2026/02/21"
$ws.Range("A855").ClearFormats()
$ws.Range("B855").Value = "土"
$ws.Range("C855").Value = 22
$ws.Range("D855").Value = 201

# New row 856: 2026/02/22 (Sun), 02:00, ranking 201
$ws.Range("A856").Value = "'2026/02/22"
$ws.Range("A856").ClearFormats()
$ws.Range("B856").Value = "日"
$ws.Range("C856").Value = 2
$ws.Range("D856").Value = 201
